$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the E2 cell (previously "added from local git")
$ws.Range("E2").ClearContents()

# Update B5 text from "Modified on local" to "Modified on master"
$ws.Range("B5").Value = "Modified on master"

# Update selection to B5
$ws.Range("B5").Select()
